# Update "想去人数" (want-to-go count) figures across the four sheets of the
# workbook to reflect a refreshed data pull (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 498
$ws1.Range("F5").Value = 1182
$ws1.Range("F7").Value = 198
$ws1.Range("F8").Value = 70
$ws1.Range("F9").Value = 779
$ws1.Range("F10").Value = 427
$ws1.Range("F11").Value = 57
$ws1.Range("F12").Value = 267
$ws1.Range("F15").Value = 10
$ws1.Range("F16").Value = 389
$ws1.Range("F17").Value = 6307
$ws1.Range("F21").Value = 7270
$ws1.Range("F24").Value = 3307
$ws1.Range("F25").Value = 435
$ws1.Range("F26").Value = 813
$ws1.Range("F27").Value = 4481
$ws1.Range("F28").Value = 336
$ws1.Range("F30").Value = 164
$ws1.Range("F31").Value = 1320
$ws1.Range("F32").Value = 126
$ws1.Range("F35").Value = 1032
$ws1.Range("F36").Value = 1372
$ws1.Range("F37").Value = 2091

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 37

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1175

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1175
$ws4.Range("F7").Value = 498
$ws4.Range("F8").Value = 1182
$ws4.Range("F10").Value = 198
$ws4.Range("F11").Value = 70
$ws4.Range("F12").Value = 779
$ws4.Range("F13").Value = 427
$ws4.Range("F14").Value = 57
$ws4.Range("F15").Value = 267
$ws4.Range("F19").Value = 10
$ws4.Range("F20").Value = 389
$ws4.Range("F21").Value = 6307
$ws4.Range("F25").Value = 7270
$ws4.Range("F28").Value = 3307
$ws4.Range("F29").Value = 435
$ws4.Range("F30").Value = 813
$ws4.Range("F31").Value = 4481
$ws4.Range("F32").Value = 336
$ws4.Range("F33").Value = 37
$ws4.Range("F35").Value = 164
$ws4.Range("F36").Value = 1320
$ws4.Range("F37").Value = 126
$ws4.Range("F40").Value = 1032
$ws4.Range("F41").Value = 1372
$ws4.Range("F43").Value = 2091
